$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 647.0454999999999
$ws.Range("I33").Value = 582.1177
$ws.Range("J33").Value = 867.8
$ws.Range("K33").Value = 582.1177
$ws.Range("L33").Value = 867.8
$ws.Range("M33").Value = -353.1177
$ws.Range("N33").Value = -1325.8
$ws.Range("H62").Value = 50140
$ws.Range("I62").Value = 50140
$ws.Range("K62").Value = 50140
$ws.Range("M62").Value = -49516
$ws.Range("H65").Value = 50140
$ws.Range("I65").Value = 50140
$ws.Range("K65").Value = 250700
$ws.Range("M65").Value = -247580
$ws.Range("H92").Value = 459.7647
$ws.Range("I92").Value = 115.28571
$ws.Range("K92").Value = 115.28571
$ws.Range("M92").Value = 1132.71429
$ws.Range("H96").Value = 932.1875
$ws.Range("J96").Value = 931
$ws.Range("L96").Value = 2793
$ws.Range("N96").Value = -5539
$ws.Range("H103").Value = 687.6667
$ws.Range("I103").Value = 206.5
$ws.Range("J103").Value = 1650
$ws.Range("K103").Value = 619.5
$ws.Range("L103").Value = 4950
$ws.Range("M103").Value = -33.5
$ws.Range("N103").Value = -6122
$ws.Range("H132").Value = 9162.817999999999
$ws.Range("I132").Value = 7579.1
$ws.Range("K132").Value = 22737.3
$ws.Range("M132").Value = -20207.3
$ws.Range("H136").Value = 148982.5
$ws.Range("J136").Value = 148982.5
$ws.Range("L136").Value = 148982.5
$ws.Range("N136").Value = -159182.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4453.6855
$ws.Range("I2").Value = 3996.3103
$ws.Range("J2").Value = 6664.3335
$ws.Range("K2").Value = 3996.3103
$ws.Range("L2").Value = 6664.3335
$ws.Range("M2").Value = -3883.3103
$ws.Range("N2").Value = -6890.3335
$ws.Range("H116").Value = 4453.6855
$ws.Range("I116").Value = 3996.3103
$ws.Range("J116").Value = 6664.3335
$ws.Range("K116").Value = 3996.3103
$ws.Range("L116").Value = 6664.3335
$ws.Range("M116").Value = -1702.3103
$ws.Range("N116").Value = -11252.3335
$ws.Range("H122").Value = 3839.3057
$ws.Range("I122").Value = 3464.68
$ws.Range("J122").Value = 4690.727
$ws.Range("K122").Value = 10394.04
$ws.Range("L122").Value = 14072.181
$ws.Range("M122").Value = -7944.039999999999
$ws.Range("N122").Value = -18972.181
$ws.Range("H132").Value = 41279.293
$ws.Range("I132").Value = 1753.625
$ws.Range("K132").Value = 5260.875
$ws.Range("M132").Value = -2730.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4453.6855
$ws.Range("I3").Value = 3996.3103
$ws.Range("J3").Value = 6664.3335
$ws.Range("K3").Value = 3996.3103
$ws.Range("L3").Value = 6664.3335
$ws.Range("M3").Value = -3882.3103
$ws.Range("N3").Value = -6892.3335
$ws.Range("H134").Value = 9795.166999999999
$ws.Range("I134").Value = 10124.417
$ws.Range("K134").Value = 30373.251
$ws.Range("M134").Value = -27838.251

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 5259
$ws.Range("I36").Value = 5259
$ws.Range("K36").Value = 5259
$ws.Range("M36").Value = -4871
$ws.Range("H40").Value = 5259
$ws.Range("I40").Value = 5259
$ws.Range("K40").Value = 5259
$ws.Range("M40").Value = -5099
$ws.Range("H58").Value = 11809.611
$ws.Range("I58").Value = 9130.214
$ws.Range("K58").Value = 9130.214
$ws.Range("M58").Value = -8927.214
$ws.Range("H86").Value = 22443.6
$ws.Range("I86").Value = 19919.857
$ws.Range("J86").Value = 28332.334
$ws.Range("K86").Value = 19919.857
$ws.Range("L86").Value = 28332.334
$ws.Range("M86").Value = -18796.857
$ws.Range("N86").Value = -30578.334
$ws.Range("H89").Value = 22443.6
$ws.Range("I89").Value = 19919.857
$ws.Range("J89").Value = 28332.334
$ws.Range("K89").Value = 99599.285
$ws.Range("L89").Value = 141661.67
$ws.Range("M89").Value = -93983.285
$ws.Range("N89").Value = -152893.67
$ws.Range("H94").Value = 3274.64
$ws.Range("I94").Value = 1252.3334
$ws.Range("J94").Value = 6308.1
$ws.Range("K94").Value = 1252.3334
$ws.Range("L94").Value = 6308.1
$ws.Range("M94").Value = -801.3334
$ws.Range("N94").Value = -7210.1
$ws.Range("H134").Value = 11839.414
$ws.Range("I134").Value = 11615.429
$ws.Range("K134").Value = 34846.287
$ws.Range("M134").Value = -32311.287
$ws.Range("H136").Value = 11809.611
$ws.Range("I136").Value = 9130.214
$ws.Range("K136").Value = 27390.642
$ws.Range("M136").Value = -24840.642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 314202.62
$ws.Range("I5").Value = 1735.2222
$ws.Range("J5").Value = 715946.4399999999
$ws.Range("K5").Value = 5205.6666
$ws.Range("L5").Value = 2147839.32
$ws.Range("M5").Value = -5093.6666
$ws.Range("N5").Value = -2148063.32
$ws.Range("H25").Value = 4759.8
$ws.Range("J25").Value = 4733
$ws.Range("L25").Value = 14199
$ws.Range("N25").Value = -14537
$ws.Range("H30").Value = 4759.8
$ws.Range("J30").Value = 4733
$ws.Range("L30").Value = 14199
$ws.Range("N30").Value = -14403
$ws.Range("H64").Value = 8301.75
$ws.Range("I64").Value = 5484.2
$ws.Range("J64").Value = 12997.667
$ws.Range("K64").Value = 16452.6
$ws.Range("L64").Value = 38993.001
$ws.Range("M64").Value = -16182.6
$ws.Range("N64").Value = -39533.001
$ws.Range("H67").Value = 8301.75
$ws.Range("I67").Value = 5484.2
$ws.Range("J67").Value = 12997.667
$ws.Range("K67").Value = 16452.6
$ws.Range("L67").Value = 38993.001
$ws.Range("M67").Value = -15516.6
$ws.Range("N67").Value = -40865.001
$ws.Range("H98").Value = 1563.8889
$ws.Range("I98").Value = 3000.5
$ws.Range("J98").Value = 1153.4286
$ws.Range("K98").Value = 9001.5
$ws.Range("L98").Value = 3460.2858
$ws.Range("M98").Value = -7503.5
$ws.Range("N98").Value = -6456.2858
$ws.Range("H122").Value = 111538.89
$ws.Range("J122").Value = 116594.07
$ws.Range("L122").Value = 1049346.63
$ws.Range("N122").Value = -1054246.63
$ws.Range("H135").Value = 314202.62
$ws.Range("I135").Value = 1735.2222
$ws.Range("J135").Value = 715946.4399999999
$ws.Range("K135").Value = 15616.9998
$ws.Range("L135").Value = 6443517.959999999
$ws.Range("M135").Value = -13081.9998
$ws.Range("N135").Value = -6448587.959999999
$ws.Range("H140").Value = 83334890
$ws.Range("I140").Value = 83334890
$ws.Range("K140").Value = 250004670
$ws.Range("M140").Value = -249999490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1249.4286
$ws.Range("J132").Value = 1780
$ws.Range("L132").Value = 5340
$ws.Range("N132").Value = -10400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5855.864
$ws.Range("I7").Value = 4711.125
$ws.Range("J7").Value = 6510
$ws.Range("K7").Value = 4711.125
$ws.Range("L7").Value = 6510
$ws.Range("M7").Value = -4599.125
$ws.Range("N7").Value = -6734
$ws.Range("H68").Value = 2046.5
$ws.Range("I68").Value = 1969.6
$ws.Range("K68").Value = 1969.6
$ws.Range("M68").Value = -1220.6
$ws.Range("H71").Value = 2046.5
$ws.Range("I71").Value = 1969.6
$ws.Range("K71").Value = 9848
$ws.Range("M71").Value = -6104
$ws.Range("H122").Value = 3639.1843
$ws.Range("I122").Value = 2783.2144
$ws.Range("K122").Value = 8349.643199999999
$ws.Range("M122").Value = -5899.643199999999
$ws.Range("H126").Value = 5855.864
$ws.Range("I126").Value = 4711.125
$ws.Range("J126").Value = 6510
$ws.Range("K126").Value = 14133.375
$ws.Range("L126").Value = 19530
$ws.Range("M126").Value = -11663.375
$ws.Range("N126").Value = -24470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 12999
$ws.Range("J43").Value = 12999
$ws.Range("L43").Value = 12999
$ws.Range("N43").Value = -13297
$ws.Range("H100").Value = 1876
$ws.Range("I100").Value = 1749
$ws.Range("J100").Value = 1926.8
$ws.Range("K100").Value = 3498
$ws.Range("L100").Value = 3853.6
$ws.Range("M100").Value = -2957
$ws.Range("N100").Value = -4935.6
$ws.Range("H122").Value = 3442
$ws.Range("I122").Value = 6000
$ws.Range("J122").Value = 3228.8333
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 9686.499899999999
$ws.Range("M122").Value = -15550
$ws.Range("N122").Value = -14586.4999
